$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save", using the same style as the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Save" values per data row, derived from the diff (rows 2..14)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 1
    6 = 0
    7 = 1
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
